$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '27.446.22'
$ws.Range("E2").Value = '  -0.50%  '
$ws.Range("D3").Value = '1.825.20'
$ws.Range("E3").Value = '  -1.97%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '332.02'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.58%  '
$ws.Range("E6").Value = '  -0.49%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4572'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.29%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3796'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.85%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.46'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07872'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.34%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9679'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.35%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.99'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.877'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.83%  '
$ws.Range("D14").Value = '1.822.93'
$ws.Range("E14").Value = '  -3.05%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.040'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.76%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.007'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.46%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '89.77'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.84%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06621'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001025'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.08'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.31%  '
$ws.Range("E21").Value = '  -0.49%  '
$ws.Range("D22").Value = '27.432.12'
$ws.Range("E22").Value = '  -0.52%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.337'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.78'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.306'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("D26").Value = '2.036.00'
$ws.Range("E26").Value = '  -2.99%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '155.39'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.63%  '
$ws.Range("E28").Value = '  -2.28%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.051'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.97%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.264'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '118.32'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9420'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.48%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09304'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.63%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.593'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.70%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.241'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.85%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.315'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05925'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.97%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02183'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.033'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.84%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.151'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.45%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5747'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1823'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.966'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.63%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.268'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '11.98'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5436'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.41%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.860'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '110.68'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.51%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06608'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.36%  '
$ws.Range("E50").Value = '  -0.62%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.040'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.23%  '
